$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value reads as a plain decimal number (e.g. "582.26")
# must be forced to stay stored as literal text (matching the source data, which
# uses inlineStr/shared-string text for ALL Price/Volume cells -- including ones that
# look numeric, like "1.00" or "27.20"). Flipping NumberFormat to "@" (Text) before
# the assignment stops Excel from auto-coercing the string into a float, then we
# restore the default "Normal" style so no stray style index is left on the cell.

$ws.Range("D2").Value = "68.024.78"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").Value = "3.256.73"
$ws.Range("E3").Value = "  -0.72%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.91%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -0.22%  "
$ws.Range("E9").Value = "  -3.42%  "
$ws.Range("E10").Value = "  -1.13%  "
$ws.Range("E11").Value = "  -3.54%  "
$ws.Range("D12").Value = "3.827.51"
$ws.Range("E12").Value = "  -0.66%  "
$ws.Range("E13").Value = "  +1.32%  "
$ws.Range("D14").Value = "68.041.25"
$ws.Range("E14").Value = "  -0.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.19"
$ws.Range("D15").Style = "Normal"
$ws.Range("E16").Value = "  -3.09%  "
$ws.Range("D17").Value = "3.255.12"
$ws.Range("E17").Value = "  -3.61%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.69"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "415.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.17%  "
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.97"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.99%  "
$ws.Range("E24").Value = "  -2.84%  "
$ws.Range("E25").Value = "  -4.29%  "
$ws.Range("E26").Value = "  -1.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.87%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("E29").Value = "  -2.68%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.53"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.42"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.15%  "
$ws.Range("E32").Value = "  -5.07%  "
$ws.Range("E33").Value = "  -5.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "164.22"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.51%  "
$ws.Range("E35").Value = "  -6.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.87"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "26.41"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.53%  "
$ws.Range("E38").Value = "  -4.46%  "
$ws.Range("E39").Value = "  -4.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.26"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.20%  "
$ws.Range("D41").Value = "2.618.78"
$ws.Range("E41").Value = "  -1.57%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0670"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.30%  "
$ws.Range("E43").Value = "  -5.75%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "332.99"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "24.05"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.40%  "
$ws.Range("E46").Value = "  -4.23%  "
$ws.Range("E47").Value = "  -2.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.975"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.18%  "
$ws.Range("E49").Value = "  -2.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "30.36"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.82%  "
